# WeeklyMidi.xlsx — refresh the week's menu (new week: 14/10 .. 18/10) and
# trim the now-unused Entree-4 / Dessert-2 columns back down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dates (row 2 header column, B2:B6) -----------------------------------
# Leading apostrophe keeps these as literal text (quote-prefixed, same as
# the existing "dd/mm"-look-alike strings) instead of being parsed as real
# dates.
$ws.Range("B2").Value = "'14/10"
$ws.Range("B3").Value = "'15/10"
$ws.Range("B4").Value = "'16/10"
$ws.Range("B5").Value = "'17/10"
$ws.Range("B6").Value = "'18/10"

# --- Monday (row 2) is now a lighter menu: only Entree 1 survives --------
$ws.Range("D2:K2").ClearContents()

# --- Tuesday (row 3) -------------------------------------------------------
$ws.Range("D3").Value = "Tomate mozarella"
$ws.Range("E3").Value = "Salade verte parisienne"
$ws.Range("F3").Value = "Feuilleté montagnard / VG"
$ws.Range("G3").Value = "Bœuf pot au feu / VG"
$ws.Range("H3").Value = "Pommes bio vapeur"
$ws.Range("I3").Value = "Filet de poisson pané MSC"
$ws.Range("J3").Value = "Légumes au pot"
$ws.Range("K3").Value = "Tartelette au citron"

# --- Wednesday (row 4) ------------------------------------------------------
$ws.Range("D4").Value = "Chicken salade"
$ws.Range("E4").Value = "Salade antillaise"
$ws.Range("F4").Value = "Soupe de légume maison"
$ws.Range("G4").Value = "Sauté de lapin à la moutarde"
$ws.Range("H4").Value = "Couscous bio"
$ws.Range("I4").Value = "Travers de porc braisé"
$ws.Range("J4").Value = "Poêlée rustique"

# --- Thursday (row 5) -------------------------------------------------------
$ws.Range("D5").Value = "Salade d'endive aux noix"
$ws.Range("E5").Value = "Crumble poireaux jambon / VG"
$ws.Range("F5").Value = "Salade athena"
$ws.Range("G5").Value = "Bœuf pot au feu / VG"
$ws.Range("H5").Value = "Lentilles bio"
$ws.Range("I5").Value = "Pavé de dinde tandoori"
$ws.Range("J5").Value = "Carottes à la crème"
$ws.Range("K5").Value = "Panna cotta fruits rouges"

# --- Friday (row 6) — also drops the Dessert-2 column (K6) -----------------
$ws.Range("D6").Value = "Betterave bio mimosa"
$ws.Range("E6").Value = "Salade exotique"
$ws.Range("F6").Value = "Salade drugstore / VG"
$ws.Range("G6").Value = "Boulettes de bœuf champignons"
$ws.Range("H6").Value = "Spaghettis"
$ws.Range("I6").Value = "Boulettes pois chiches bio"
$ws.Range("J6").Value = "Épinard béchamel"
$ws.Range("K6").ClearContents()

# --- Widen the Dessert 1 column (K) now that its entries are longer -------
$ws.Columns.Item(11).ColumnWidth = 22.43
